# Apply the edit described by the diff:
#  - rename the sheet
#  - insert a "new" day of data at the top of the table (row 2), pushing the
#    data that used to occupy rows 2..157 down into rows 3..158
#  - the row that used to be at row 158 (date 44070) is relocated to a brand
#    new row 396 appended after the existing last row (395); rows 159..395
#    are left untouched
#  - fix up the worksheet dimension / style of the newly created row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet
$ws.Name = "overview_2021-01-31"

# 2. Remember the data currently sitting in row 158 -- it needs to move to
#    the new row 396 at the bottom of the sheet.
$savedDate = $ws.Cells.Item(158, 1).Value2()
$savedB    = $ws.Cells.Item(158, 2).Value()
$savedC    = $ws.Cells.Item(158, 3).Value()
$savedD    = $ws.Cells.Item(158, 4).Value()
$savedE    = $ws.Cells.Item(158, 5).Value()
$savedF    = $ws.Cells.Item(158, 6).Value()

# 3. Shift rows 2..157 down into rows 3..158 (working from the bottom up so
#    we never overwrite a source row before it has been read).
for ($k = 157; $k -ge 2; $k--) {
    $srcDate = $ws.Cells.Item($k, 1).Value2()
    $srcB    = $ws.Cells.Item($k, 2).Value()
    $srcC    = $ws.Cells.Item($k, 3).Value()
    $srcD    = $ws.Cells.Item($k, 4).Value()
    $srcE    = $ws.Cells.Item($k, 5).Value()
    $srcF    = $ws.Cells.Item($k, 6).Value()

    $dst = $k + 1
    $ws.Cells.Item($dst, 1).Value = $srcDate
    $ws.Cells.Item($dst, 2).Value = $srcB
    $ws.Cells.Item($dst, 3).Value = $srcC
    $ws.Cells.Item($dst, 4).Value = $srcD
    $ws.Cells.Item($dst, 5).Value = $srcE
    $ws.Cells.Item($dst, 6).Value = $srcF
}

# 4. Write the brand new row of data into row 2.
$ws.Cells.Item(2, 1).Value = 44227
$ws.Cells.Item(2, 2).Value = "overview"
$ws.Cells.Item(2, 3).Value = "K02000001"
$ws.Cells.Item(2, 4).Value = "United Kingdom"
$ws.Cells.Item(2, 5).Value = 21088
$ws.Cells.Item(2, 6).Value = 587

# 5. Append the data that used to be in row 158 as new row 396, copying the
#    date style (s="1") from the existing last row (395) first.
$ws.Range("A395").Copy()
$ws.Range("A396").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(396, 1).Value = $savedDate
$ws.Cells.Item(396, 2).Value = $savedB
$ws.Cells.Item(396, 3).Value = $savedC
$ws.Cells.Item(396, 4).Value = $savedD
$ws.Cells.Item(396, 5).Value = $savedE
$ws.Cells.Item(396, 6).Value = $savedF
